# Addition to last commit:
# Shift the recorded "Start of the project" dates (and the related
# Notes timestamp) forward by 3 days:
#   43785.0              -> 43788.0               (2019-11-16 -> 2019-11-19)
#   43785.334085648145   -> 43788.81340277778      (2019-11-16 08:01:05 -> 2019-11-19 19:31:18)

$wb = $excel.ActiveWorkbook

# Projects sheet: C2 = "Start of the project" for the testovaci_projekt row.
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Range("C2").Value2 = 43788.0

# Tasks sheet: D2/D3 = "Start of the project" for both testTask rows.
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("D2").Value2 = 43788.0
$wsTasks.Range("D3").Value2 = 43788.0

# Notes sheet: B2 = Timestamp for the testovaci text note.
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("B2").Value2 = 43788.81340277778
